$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($count))
$ws.Name = "Bugs"

$ws.Range("A1").Value = "File"
$ws.Range("B1").Value = "Function"
$ws.Range("B2").Value = "isnumeric()"
$ws.Range("C1").Value = "Frequency"
$ws.Range("A2").Value = "Everywhere"
$ws.Range("B3").Value = "GetParam()"
$ws.Range("C2").Value = "High"
$ws.Range("C3").Value = "Low"
$ws.Range("A3").Value = "Machine"

$ws.Columns.Item(1).ColumnWidth = 9.833333333333334
$ws.Columns.Item(2).ColumnWidth = 9.333333333333334
$ws.Columns.Item(3).ColumnWidth = 9.333333333333334

$ws.Range("C4").Select() | Out-Null
